{"js": "// The original page had two pictures (inline drawings) placed in their own\n// paragraphs, and the text \"paper clips'Gn\" used a straight apostrophe.\n// The target edit removes the two picture paragraphs (pipeline now handles\n// images separately) and normalizes the apostrophe to a curly quote.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Paragraphs that consist solely of an inline picture have empty text\n// (the drawing itself isn't reflected in Range.text). Collect and delete\n// those paragraph objects (this removes the drawing along with the\n// paragraph mark, matching the diff which drops the whole <w:p>).\nconst emptyParagraphs = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\") {\n    emptyParagraphs.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of emptyParagraphs) {\n  p.delete();\n}\nawait context.sync();\n\n// Fix the mangled apostrophe: \"clips'Gn\" -> \"clips\\u2019Gn\" (U+2019).\nconst searchResults = body.search(\"paper clips'Gn\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"paper clips\\u2019Gn\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The original page had two pictures (inline drawings), each sitting alone\n# inside its own paragraph, and the text \"paper clips'Gn\" used a straight\n# apostrophe. The target edit removes the two picture-only paragraphs\n# (image handling moved elsewhere in the pipeline) and normalizes the\n# apostrophe to a curly right single quotation mark (U+2019).\n\n$d = $word.ActiveDocument\n\n# Paragraphs that hold nothing but an inline picture have a Range.Text\n# consisting of just the paragraph mark (chr 13); walk backwards so the\n# indices of paragraphs still to be inspected stay valid as we delete.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $textOnly = $p.Range.Text.Replace([char]13, \"\")\n    if ($textOnly -eq \"\") {\n        $p.Range.Delete()\n    }\n}\n\n# Fix the mangled apostrophe: \"clips'Gn\" -> \"clips\" + U+2019 + \"Gn\".\n$searchText = \"paper clips\" + [char]0x27 + \"Gn\"\n$replaceText = \"paper clips\" + [char]0x2019 + \"Gn\"\n\n$find = $d.Content.Find\n$find.Text = $searchText\n$find.Replacement.Text = $replaceText\n$find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
